# Render site for cm012 and hw07
# - Row 14 (cm012) topic renamed: "Distributed computing" -> "Distributed learning"
# - Row 14 link_it flag flipped from FALSE to TRUE (slides now linked)
# - Active cell/selection moved from C14 to C15 after the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D14").Value = "Distributed learning"
$ws.Range("C14").Value = $true

$ws.Range("C15").Select() | Out-Null
